$d = $word.ActiveDocument

# Build the old/new strings with explicit characters to avoid any
# encoding ambiguity (straight apostrophe, middle dot, accented vowels, cedilla).
$old = "Dates de la campanya constel" + [char]0xB7 + "laci" + [char]0xF3 + " d'Ori" + [char]0xF3 + " 2022: 16-25 de gener, 14-23 de febrer, 14-24 de mar" + [char]0xE7
$new = "Dates de la campanya 2022 en qu" + [char]0xE8 + " usem la constel" + [char]0xB7 + "laci" + [char]0xF3 + ", constel" + [char]0xB7 + "laci" + [char]0xF3 + " d'Ori" + [char]0xF3 + " 16-25 de gener, 14-23 de febrer, 14-24 de mar" + [char]0xE7

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = $old
$find.Forward = $true
$find.Wrap = 1
$find.MatchCase = $true
$find.MatchWholeWord = $false
$find.MatchWildcards = $false

# Repeatedly find each occurrence and overwrite the located Range's .Text
# directly (rather than using Find.Execute's Replace parameter), because
# the Replace pathway runs AutoCorrect/AutoFormat-as-you-type and mangles
# the straight apostrophe into a curly one. Setting Range.Text directly
# inserts the literal characters untouched.
$keepGoing = $true
while ($keepGoing) {
    $found = $find.Execute()
    if ($found) {
        $rng = $find.Parent
        $rng.Text = $new
    } else {
        $keepGoing = $false
    }
}
